$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (14-26) appended to the smell-test log.
# Row -> (Timestamp, Refrigerator Temp, Refrigerated Milk Smell)
$rowData = @{
    14 = @("12-06-2024 08:46AM", 38,   "Slightly Sour")
    15 = @("12-06-2024 09:55PM", 38.5, "Slightly Sour")
    16 = @("12-07-2024 09:22AM", 38,   "Slightly Sour")
    17 = @("12-07-2024 08:05PM", 38,   "Slightly Sour")
    18 = @("12-08-2024 09:31AM", 38.5, "Slightly Sour")
    19 = @("12-08-2024 07:12AM", 37.5, "Slightly Sour")
    20 = @("12-09-2024 06:20AM", 37,   "Slightly Sour")
    21 = @("12-09-2024 09:07PM", 37,   "Slightly Sour")
    22 = @("12-10-2024 06:18AM", 38,   "Slightly Sour")
    23 = @("12-10-2024 06:20PM", 38,   "Noticeably Sour")
    24 = @("12-11-2024 06:17AM", 38.5, "Noticeably Sour")
    25 = @("12-11-2024 08:55PM", 37,   "Noticeably Sour")
    26 = @("12-12-2024 10:46AM", 38.5, "Noticeably Sour")
}

# Timestamps (column A) were originally typed out of strict row order -- replicate
# the exact entry order so the workbook's shared-string table lines up.
$timestampEntryOrder = @(14, 15, 16, 17, 18, 19, 20, 22, 23, 21, 24, 25, 26)

foreach ($row in $timestampEntryOrder) {
    $ws.Cells.Item($row, 1).Value = $rowData[$row][0]
}

foreach ($row in 14..26) {
    $ws.Cells.Item($row, 2).Value = $rowData[$row][1]
    $ws.Cells.Item($row, 3).Value = $rowData[$row][2]
}

$ws.Range("C26").Select()
